$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-9 and add new rows 10-13 for the Inhbb-Acvr1b LR pair table
$stringData = @{
  2 = @("ECs","Inhbb","Acvr1b","ECs")
  3 = @("ECs","Inhbb","Acvr1b","FAPs")
  4 = @("ECs","Inhbb","Acvr1b","M2")
  5 = @("ECs","Inhbb","Acvr1b","sCs")
  6 = @("FAPs","Inhbb","Acvr1b","ECs")
  7 = @("FAPs","Inhbb","Acvr1b","FAPs")
  8 = @("FAPs","Inhbb","Acvr1b","M2")
  9 = @("FAPs","Inhbb","Acvr1b","sCs")
  10 = @("sCs","Inhbb","Acvr1b","ECs")
  11 = @("sCs","Inhbb","Acvr1b","FAPs")
  12 = @("sCs","Inhbb","Acvr1b","M2")
  13 = @("sCs","Inhbb","Acvr1b","sCs")
}

$numData = @{
  2 = @(3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,2.548026,7.644078,0.1809961777378536,0.1809961777378536,6.434157196662,57.90741476995801,0.1335424507626665,0.1335424507626666)
  3 = @(3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,6.757657999999999,20.272974,0.4800226796977849,0.4800226796977848,17.06412487677933,153.577123891014,0.3541699381152074,0.3541699381152075)
  4 = @(3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,3.246115666666667,9.738347000000001,0.2305841966140185,0.2305841966140185,8.196940878107446,73.77246790296701,0.1701294420016726,0.1701294420016726)
  5 = @(3,1,2.525153666666667,7.575461000000001,0.7378191762484796,0.7378191762484798,3,1,1.525989333333333,4.577968,0.108396945950343,0.108396945950343,3.853357560360889,34.680218043248,0.07997734536893306,0.07997734536893307)
  6 = @(2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,2.548026,7.644078,0.1809961777378536,0.1809961777378536,0.9528997220340002,8.576097498306002,0.01977765856847605,0.01977765856847606)
  7 = @(2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,6.757657999999999,20.272974,0.4800226796977849,0.4800226796977848,2.527199655655333,22.744796900898,0.05245262514846032,0.05245262514846032)
  8 = @(2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,3.246115666666667,9.738347000000001,0.2305841966140185,0.2305841966140185,1.213968270518778,10.925714434669,0.02519619789166766,0.02519619789166766)
  9 = @(2,0.6666666666666666,0.3739756666666667,1.121927,0.1092711394000877,0.1092711394000878,3,1,1.525989333333333,4.577968,0.108396945950343,0.108396945950343,0.5706828782595557,5.136145904336001,0.01184465779148371,0.01184465779148371)
  10 = @(3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,2.548026,7.644078,0.1809961777378536,0.1809961777378536,1.33344995316,12.00104957844,0.02767606840671098,0.02767606840671099)
  11 = @(3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,6.757657999999999,20.272974,0.4800226796977849,0.4800226796977848,3.536462635613332,31.82816372052,0.07340011643411712,0.07340011643411712)
  12 = @(3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,3.246115666666667,9.738347000000001,0.2305841966140185,0.2305841966140185,1.698778891451111,15.28901002306,0.03525855672067824,0.03525855672067825)
  13 = @(3,1,0.5233266666666666,1.56998,0.1529096843514326,0.1529096843514326,3,1,1.525989333333333,4.577968,0.108396945950343,0.108396945950343,0.7985909111822221,7.18731820064,0.01657494278992625,0.01657494278992625)
}

foreach ($r in 2..13) {
  $sd = $stringData[$r]
  $ws.Cells.Item($r, 1).Value = $sd[0]
  $ws.Cells.Item($r, 2).Value = $sd[1]
  $ws.Cells.Item($r, 3).Value = $sd[2]
  $ws.Cells.Item($r, 4).Value = $sd[3]
  $nd = $numData[$r]
  for ($i = 0; $i -lt $nd.Length; $i++) {
    $ws.Cells.Item($r, 5 + $i).Value = $nd[$i]
  }
}
